$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Checklist")

# Insert a new row at position 49 (shifts rows 49-132 down to 50-133).
$ws.Range("A49").EntireRow.Insert()

# New row height (66pt, matches other multi-line rows in the checklist).
$ws.Rows.Item(49).RowHeight = 66

# Copy the formatting from the row directly above (WSTG-ATHN-10, row 48)
# onto the freshly inserted row so the new entry matches the existing
# checklist look (borders/fonts/alignment for columns A-F).
$fmtSrc = $ws.Range("A48:F48")
$fmtDst = $ws.Range("A49:F49")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new WSTG-AUTH-11 row content.
$ws.Range("A49").Value = ""
$ws.Range("B49").Value = "WSTG-AUTH-11"
$ws.Range("C49").Formula = '=HYPERLINK("https://owasp.org/www-project-web-security-testing-guide/latest/4-Web_Application_Security_Testing/04-Authentication_Testing/11-Testing_Multi-Factor_Authentication", "Testing Multi-Factor Authentication (MFA)")'
$ws.Range("D49").Value = "- Identify the type of MFA used by the application.
- Determine whether the MFA implementation is robust and secure.
- Attempt to bypass the MFA."
$ws.Range("E49").Value = "Not Started"
$ws.Range("F49").Value = ""

# Extend the conditional formatting range that covers the data rows so it
# includes the new last row (133).
$ws.Range("B4:F133").FormatConditions.Delete()
$src = $ws.Range("B4:F4")
$src.Copy()
$dstCf = $ws.Range("B4:F133")
$dstCf.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Make sure the new row's Status cell (E49) has the same data-validation
# dropdown ("Not Started, Pass, Issues, N/A") as the other status cells.
$ws.Range("E49").Validation.Add(3, 1, 1, "Not Started,Pass,Issues,N/A")
$ws.Range("E49").Validation.IgnoreBlank = $true
$ws.Range("E49").Validation.InCellDropdown = $true
